# Word namespace used throughout for InsertXML fragments.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart .../><w:bookmarkEnd .../> pair (the
#    "_GoBack" bookmark) that currently sits between the "Entity framework"
#    run and the " (Not completed)" run, inside the ORM bullet paragraph.
# ---------------------------------------------------------------------------
$ormPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Entity framework*") {
        $ormPara = $candidate
        break
    }
}

if ($ormPara -ne $null) {
    # NB: the rsid attributes below are copied verbatim from the source
    # document so this whole-paragraph replace is a byte-for-byte no-op
    # except for the removal of the bookmark pair.
    $xml = "<w:p $wns w:rsidR=`"008B2BD6`" w:rsidRDefault=`"008B2BD6`" w:rsidP=`"008B2BD6`"><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"5`"/></w:numPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"FFFFFF`"/><w:spacing w:before=`"60`" w:after=`"100`" w:afterAutospacing=`"1`" w:line=`"240`" w:lineRule=`"auto`"/><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:color w:val=`"24292E`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:color w:val=`"24292E`"/></w:rPr><w:t xml:space=`"preserve`">ORM implementation of classes with </w:t></w:r><w:r w:rsidR=`"007E7FB1`"><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:color w:val=`"24292E`"/></w:rPr><w:t>Entity framework</w:t></w:r><w:r w:rsidR=`"002F5E06`"><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:color w:val=`"24292E`"/></w:rPr><w:t xml:space=`"preserve`"> (Not completed)</w:t></w:r></w:p>"
    $ormPara.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) Replace the trailing empty paragraph at the end of the document with a
#    new "Note:" paragraph describing where the running application can be
#    reached. The relocated "_GoBack" bookmark is re-created inside this new
#    paragraph, right after the bold "Note: " label.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$noteXml = "<w:p $wns><w:pPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:b/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>Note:</w:t></w:r>" +
    "<w:r><w:rPr><w:b/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "<w:r><w:rPr><w:b/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:br/></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Application runs with all features on </w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>localhost:4200</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> and </w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>localhost:44355</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>/api for front and backends.</w:t></w:r>" +
    "</w:p>"

$lastPara.Range.InsertXML($noteXml)

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
